$d = $word.ActiveDocument

# The final paragraph of the document holds the _GoBack bookmark. We need to:
#   1. insert a brand-new empty paragraph right before it, and
#   2. prepend three runs ("Changes made by ", "manvir" wrapped in spell-check
#      markers, and a trailing space) to that bookmark paragraph, leaving the
#      bookmarkStart/bookmarkEnd in place at the end.
#
# Range.InsertXML splices a WordprocessingML fragment in place of the
# paragraph that "owns" the range it's called on, and (because that paragraph
# is the very last one in the body) leaves one extra empty trailing paragraph
# behind to keep the body's final paragraph mark intact. So we build the
# fragment to contain both the new empty spacer paragraph *and* the updated
# bookmark paragraph, then fold the stray trailing paragraph mark back out
# afterwards.

$lastPara = $d.Paragraphs.Last
$insertionRange = $lastPara.Range
$insertionRange.Collapse(1)

$fragment = @'
<?xml version="1.0"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p/>
          <w:p>
            <w:r>
              <w:t xml:space="preserve">Changes made by </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>manvir</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$insertionRange.InsertXML($fragment)

# InsertXML above leaves a stray empty paragraph after the bookmark paragraph
# (the remnant of the original final paragraph mark). Merge it away by
# deleting the paragraph mark that now separates the bookmark paragraph from
# that stray trailing paragraph.
$bookmarkPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$mergeRange = $d.Range($bookmarkPara.Range.End - 1, $bookmarkPara.Range.End)
$mergeRange.Delete()
